$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the trailing numeric placeholder values (0) from F48:G49 so these
# footer rows (e.g. "Back to index" / "Please click to email us your opinion:")
# no longer carry stray 0 values in the Gap_Growth_%/Ratio_Change_% columns.
$ws.Range("F48:G49").ClearContents()
